# Applies the "Vergelijking Specificaties" table edit:
#   - widen/narrow the 4 table columns (1685/2563/2876/2523 -> 2116/2419/2712/2400 dxa)
#   - relabel the "Video Encoding" row to "Video Encoding/Decoding"
#     (folding the old, separate "Video Decoding" row's heading into it)
#   - turn the former "Video Decoding" row into a new "Power (V/A)" spec row
#     with values 12V/5A, 12V/4.8A, 5V/4A
#   - append a new, still-empty row at the bottom of the table

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Resize the four columns (dxa -> points is dxa/20) ---------------
$t.Columns.Item(1).Width = 2116 / 20
$t.Columns.Item(2).Width = 2419 / 20
$t.Columns.Item(3).Width = 2712 / 20
$t.Columns.Item(4).Width = 2400 / 20

# --- 2. "Video Encoding" -> "Video Encoding/Decoding" --------------------
$found = $d.Content.Find.Execute("Video Encoding", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Video Encoding/Decoding", 2)

# --- 3. Turn the old "Video Decoding" row into the new "Power (V/A)" row -
$row = $t.Rows.Item(10)
$row.Cells.Item(1).Range.Text = "Power (V/A)"
$row.Cells.Item(2).Range.Text = "12V/5A"
$row.Cells.Item(3).Range.Text = "12V/4.8A"
$row.Cells.Item(4).Range.Text = "5V/4A"

# --- 4. Append a fresh, empty row at the end of the table -----------------
$newRow = $t.Rows.Add()
